$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.93462456970575
$ws.Range("C2").Value = 5.280088797411882
$ws.Range("D2").Value = 13.61372221487566
$ws.Range("E2").Value = 13.86365087178498
$ws.Range("G2").Value = 3.711066082556555
$ws.Range("I2").Value = 30.10502409093299
$ws.Range("J2").Value = 8.353572073523754
$ws.Range("K2").Value = 12.98561771456335
$ws.Range("L2").Value = 12.65939286758939
$ws.Range("N2").Value = 21.76002348818933
$ws.Range("O2").Value = 31.66240120664786

$ws.Range("B3").Value = 15.76959921898604
$ws.Range("C3").Value = 5.18697315462554
$ws.Range("D3").Value = 13.60705870072544
$ws.Range("E3").Value = 13.88313094264151
$ws.Range("G3").Value = 3.713048051616156
$ws.Range("I3").Value = 30.19097428770155
$ws.Range("J3").Value = 8.362569526578678
$ws.Range("K3").Value = 12.87267849026182
$ws.Range("L3").Value = 12.66414272467651
$ws.Range("N3").Value = 21.82182743811731
$ws.Range("O3").Value = 31.74084142498133

$ws.Range("B4").Value = 15.67039747458705
$ws.Range("C4").Value = 5.12812818411575
$ws.Range("D4").Value = 13.60535356381548
$ws.Range("E4").Value = 13.89686848301197
$ws.Range("G4").Value = 3.714330535253257
$ws.Range("I4").Value = 30.24832114261384
$ws.Range("J4").Value = 8.368406131399338
$ws.Range("K4").Value = 12.80510504923499
$ws.Range("L4").Value = 12.66878799187089
$ws.Range("N4").Value = 21.86155322336774
$ws.Range("O4").Value = 31.79428019868575

$ws.Range("B5").Value = 15.63054865409754
$ws.Range("C5").Value = 5.103743146874624
$ws.Range("D5").Value = 13.60526065224379
$ws.Range("E5").Value = 13.90291387127541
$ws.Range("G5").Value = 3.714869691400505
$ws.Range("I5").Value = 30.27284017480415
$ws.Range("J5").Value = 8.370863315450205
$ws.Range("K5").Value = 12.77804010356867
$ws.Range("L5").Value = 12.67111664690625
$ws.Range("N5").Value = 21.87819017832594
$ws.Range("O5").Value = 31.81738269274426

$ws.Range("B6").Value = 15.62396783938127
$ws.Range("C6").Value = 5.099670009274392
$ws.Range("D6").Value = 13.60528162982889
$ws.Range("E6").Value = 13.90394472762228
$ws.Range("G6").Value = 3.714960217974265
$ws.Range("I6").Value = 30.27698097531408
$ws.Range("J6").Value = 8.3712760903268
$ws.Range("K6").Value = 12.77357523567883
$ws.Range("L6").Value = 12.67152965709368
$ws.Range("N6").Value = 21.88097985116003
$ws.Range("O6").Value = 31.82129888676132

$ws.Range("B7").Value = 15.66985766910671
$ws.Range("C7").Value = 5.12780093930836
$ws.Range("D7").Value = 13.60534987123738
$ws.Range("E7").Value = 13.89694820183334
$ws.Range("G7").Value = 3.714337739489521
$ws.Range("I7").Value = 30.24864715997947
$ws.Range("J7").Value = 8.36843895079803
$ws.Range("K7").Value = 12.8047380976079
$ws.Range("L7").Value = 12.66881763174827
$ws.Range("N7").Value = 21.8617757779823
$ws.Range("O7").Value = 31.79458640041649

$ws.Range("B8").Value = 15.87730867652874
$ws.Range("C8").Value = 5.248338306956781
$ws.Range("D8").Value = 13.61093058367573
$ws.Range("E8").Value = 13.86999913686068
$ws.Range("G8").Value = 3.711735892283512
$ws.Range("I8").Value = 30.13371045087555
$ws.Range("J8").Value = 8.356609752035402
$ws.Range("K8").Value = 12.94632455822728
$ws.Range("L8").Value = 12.66067239320373
$ws.Range("N8").Value = 21.78096527414563
$ws.Range("O8").Value = 31.68835176597281

$ws.Range("B9").Value = 16.29906716232697
$ws.Range("C9").Value = 5.470821985966389
$ws.Range("D9").Value = 13.64071386311481
$ws.Range("E9").Value = 13.83122810336316
$ws.Range("G9").Value = 3.70715141643069
$ws.Range("I9").Value = 29.94461470672587
$ws.Range("J9").Value = 8.335878819688501
$ws.Range("K9").Value = 13.23684323034009
$ws.Range("L9").Value = 12.65837256539859
$ws.Range("N9").Value = 21.63654242800096
$ws.Range("O9").Value = 31.52193138004179

$ws.Range("B10").Value = 16.61544747854923
$ws.Range("C10").Value = 5.625027991884597
$ws.Range("D10").Value = 13.6739371170584
$ws.Range("E10").Value = 13.81129406938124
$ws.Range("G10").Value = 3.704095550603129
$ws.Range("I10").Value = 29.82782470256025
$ws.Range("J10").Value = 8.322136514173573
$ws.Range("K10").Value = 13.45649893834834
$ws.Range("L10").Value = 12.66495361832119
$ws.Range("N10").Value = 21.53891160936789
$ws.Range("O10").Value = 31.42526231826379

$ws.Range("B11").Value = 16.76023492442525
$ws.Range("C11").Value = 5.693019911025234
$ws.Range("D11").Value = 13.69147698859853
$ws.Range("E11").Value = 13.80407499483391
$ws.Range("G11").Value = 3.702772484775671
$ws.Range("I11").Value = 29.77950354834237
$ws.Range("J11").Value = 8.316204927502675
$ws.Range("K11").Value = 13.55741986228975
$ws.Range("L11").Value = 12.66972730294304
$ws.Range("N11").Value = 21.49631900903817
$ws.Range("O11").Value = 31.38685341796693

$ws.Range("B12").Value = 16.81514007253747
$ws.Range("C12").Value = 5.718444779020479
$ws.Range("D12").Value = 13.69846412822311
$ws.Range("E12").Value = 13.80160648700689
$ws.Range("G12").Value = 3.702281064923874
$ws.Range("I12").Value = 29.76189708672057
$ws.Range("J12").Value = 8.314004548390148
$ws.Range("K12").Value = 13.5957495605784
$ws.Range("L12").Value = 12.67178932627472
$ws.Range("N12").Value = 21.48045067570119
$ws.Range("O12").Value = 31.3731100100103

$ws.Range("B13").Value = 16.80331257185509
$ws.Range("C13").Value = 5.71298359140306
$ws.Range("D13").Value = 13.69694403412874
$ws.Range("E13").Value = 13.80212634109555
$ws.Range("G13").Value = 3.702386474921991
$ws.Range("I13").Value = 29.76565818248621
$ws.Range("J13").Value = 8.31447640648592
$ws.Range("K13").Value = 13.58749002617045
$ws.Range("L13").Value = 12.67133394827147
$ws.Range("N13").Value = 21.48385663532284
$ws.Range("O13").Value = 31.3760342537826

$ws.Range("B14").Value = 16.76475072242849
$ws.Range("C14").Value = 5.695118153493815
$ws.Range("D14").Value = 13.69204492962082
$ws.Range("E14").Value = 13.80386659836655
$ws.Range("G14").Value = 3.702731863309435
$ws.Range("I14").Value = 29.77804118880246
$ws.Range("J14").Value = 8.316022984610886
$ws.Range("K14").Value = 13.5605711762692
$ws.Range("L14").Value = 12.66989185936252
$ws.Range("N14").Value = 21.49500829638568
$ws.Range("O14").Value = 31.38570667809467

$ws.Range("B15").Value = 16.74113915488745
$ws.Range("C15").Value = 5.684132754530022
$ws.Range("D15").Value = 13.68908892150937
$ws.Range("E15").Value = 13.80496707155653
$ws.Range("G15").Value = 3.702944672021871
$ws.Range("I15").Value = 29.78571623152056
$ws.Range("J15").Value = 8.316976264569693
$ws.Range("K15").Value = 13.54409643499997
$ws.Range("L15").Value = 12.66904160966672
$ws.Range("N15").Value = 21.50187290870104
$ws.Range("O15").Value = 31.3917356792644

$ws.Range("B16").Value = 16.60599864188998
$ws.Range("C16").Value = 5.620540109118602
$ws.Range("D16").Value = 13.67283935816797
$ws.Range("E16").Value = 13.81180300238286
$ws.Range("G16").Value = 3.704183361507673
$ws.Range("I16").Value = 29.83107937511987
$ws.Range("J16").Value = 8.322530575572094
$ws.Range("K16").Value = 13.4499210061971
$ws.Range("L16").Value = 12.66467732801658
$ws.Range("N16").Value = 21.54173166130712
$ws.Range("O16").Value = 31.42788450835828

$ws.Range("B17").Value = 16.52328210787208
$ws.Range("C17").Value = 5.580966797271345
$ws.Range("D17").Value = 13.66348968056657
$ws.Range("E17").Value = 13.81646973127547
$ws.Range("G17").Value = 3.704960400160537
$ws.Range("I17").Value = 29.86013975632687
$ws.Range("J17").Value = 8.326019737787263
$ws.Range("K17").Value = 13.39238148210542
$ws.Range("L17").Value = 12.66245475341928
$ws.Range("N17").Value = 21.56664905914334
$ws.Range("O17").Value = 31.45148690832341

$ws.Range("B18").Value = 16.47578990025054
$ws.Range("C18").Value = 5.558002900616363
$ws.Range("D18").Value = 13.65834062631344
$ws.Range("E18").Value = 13.81932799330822
$ws.Range("G18").Value = 3.705413647306889
$ws.Range("I18").Value = 29.87730697095618
$ws.Range("J18").Value = 8.328056730566919
$ws.Range("K18").Value = 13.35938206991641
$ws.Range("L18").Value = 12.66134402766567
$ws.Range("N18").Value = 21.58115227787161
$ws.Range("O18").Value = 31.46558627598909

$ws.Range("B19").Value = 16.45972574957942
$ws.Range("C19").Value = 5.550193340002991
$ws.Range("D19").Value = 13.65663662205822
$ws.Range("E19").Value = 13.82032567065192
$ws.Range("G19").Value = 3.705568195063018
$ws.Range("I19").Value = 29.88319719010246
$ws.Range("J19").Value = 8.328751600742025
$ws.Range("K19").Value = 13.34822641817334
$ws.Range("L19").Value = 12.66099679467977
$ws.Range("N19").Value = 21.58609229102431
$ws.Range("O19").Value = 31.47045003287089

$ws.Range("B20").Value = 16.53207905670786
$ws.Range("C20").Value = 5.585200478220933
$ws.Range("D20").Value = 13.66446133322208
$ws.Range("E20").Value = 13.81595493788302
$ws.Range("G20").Value = 3.704877029826272
$ws.Range("I20").Value = 29.85699939876709
$ws.Range("J20").Value = 8.325645194849193
$ws.Range("K20").Value = 13.39849696843376
$ws.Range("L20").Value = 12.66267401146192
$ws.Range("N20").Value = 21.56397882986531
$ws.Range("O20").Value = 31.44892016572605

$ws.Range("B21").Value = 16.7760755468081
$ws.Range("C21").Value = 5.7003745017845
$ws.Range("D21").Value = 13.69347457895436
$ws.Range("E21").Value = 13.80334825095959
$ws.Range("G21").Value = 3.702630154295949
$ws.Range("I21").Value = 29.77438522346367
$ws.Range("J21").Value = 8.315567476233031
$ws.Range("K21").Value = 13.56847506668637
$ws.Range("L21").Value = 12.67030854684839
$ws.Range("N21").Value = 21.4917257209976
$ws.Range("O21").Value = 31.38284390245489

$ws.Range("B22").Value = 16.93597046917232
$ws.Range("C22").Value = 5.773764354555144
$ws.Range("D22").Value = 13.71444637500313
$ws.Range("E22").Value = 13.79665447556553
$ws.Range("G22").Value = 3.701217605202125
$ws.Range("I22").Value = 29.724424089681
$ws.Range("J22").Value = 8.309247872286443
$ws.Range("K22").Value = 13.6802109931041
$ws.Range("L22").Value = 12.67677976189136
$ws.Range("N22").Value = 21.44602234327217
$ws.Range("O22").Value = 31.34432957602853

$ws.Range("B23").Value = 16.85060765353452
$ws.Range("C23").Value = 5.734770871437123
$ws.Range("D23").Value = 13.70307073976918
$ws.Range("E23").Value = 13.800085905108
$ws.Range("G23").Value = 3.701966408566471
$ws.Range("I23").Value = 29.75072021948512
$ws.Range("J23").Value = 8.312596423821649
$ws.Range("K23").Value = 13.62052639596904
$ws.Range("L23").Value = 12.67319095317551
$ws.Range("N23").Value = 21.4702765656118
$ws.Range("O23").Value = 31.36445783693558

$ws.Range("B24").Value = 16.52810175745483
$ws.Range("C24").Value = 5.583287092598612
$ws.Range("D24").Value = 13.66402134419842
$ws.Range("E24").Value = 13.81618712989749
$ws.Range("G24").Value = 3.704914701249916
$ws.Range("I24").Value = 29.85841772151823
$ws.Range("J24").Value = 8.325814429050531
$ws.Range("K24").Value = 13.39573190296922
$ws.Range("L24").Value = 12.66257436440164
$ws.Range("N24").Value = 21.56518548634595
$ws.Range("O24").Value = 31.45007893907187

$ws.Range("B25").Value = 16.18363700018467
$ws.Range("C25").Value = 5.412205201079789
$ws.Range("D25").Value = 13.63065336913754
$ws.Range("E25").Value = 13.84021266180868
$ws.Range("G25").Value = 3.70833655151411
$ws.Range("I25").Value = 29.99188316217996
$ws.Range("J25").Value = 8.341224599031893
$ws.Range("K25").Value = 13.15703948236856
$ws.Range("L25").Value = 12.65753714096026
$ws.Range("N25").Value = 21.67411784793859
$ws.Range("O25").Value = 31.5624603262518
